$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Material")

# Row 8 and Row 9 (columns B..H) had their contents swapped.
# Column A ("Sheet") and column G (empty) were already identical in both
# rows, so only B..F and H actually change value, but we swap the whole
# B:H range for correctness.

$row8 = @($ws.Range("B8").Value2, $ws.Range("C8").Value2, $ws.Range("D8").Value2, $ws.Range("E8").Value2, $ws.Range("F8").Value2, $ws.Range("G8").Value2, $ws.Range("H8").Value2)
$row9 = @($ws.Range("B9").Value2, $ws.Range("C9").Value2, $ws.Range("D9").Value2, $ws.Range("E9").Value2, $ws.Range("F9").Value2, $ws.Range("G9").Value2, $ws.Range("H9").Value2)

$cols = @("B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range("$col" + "8").Value2 = $row9[$i]
    $ws.Range("$col" + "9").Value2 = $row8[$i]
}
